# Locate the paragraph that currently ends the "signal doesn't reach" bullet
# list item (last bullet before the trailing bookmark) and extend it with two
# additional sentences, then add a brand-new bullet paragraph after it with
# matching list formatting.

$d = $word.ActiveDocument

$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*signal doesn’t reach.*") {
        $targetIndex = $i
        $found = $true
        break
    }
}

if (-not $found) {
    throw "Could not find target paragraph 'The signal doesn't reach.'"
}

# --- Append the two new sentences to the end of the paragraph's text,
#     inserted just before the paragraph mark so existing bookmark/formatting
#     at the end of the paragraph is preserved. ---
$para = $d.Paragraphs.Item($targetIndex)
$insertionPoint = $para.Range.End - 1
$ins = $d.Range($insertionPoint, $insertionPoint)
$ins.InsertAfter(" This Bluetooth module is rated for 10m and drywall or other mate")

$para = $d.Paragraphs.Item($targetIndex)
$insertionPoint = $para.Range.End - 1
$ins = $d.Range($insertionPoint, $insertionPoint)
$ins.InsertAfter("rial shouldn’t stop the signal. If you think you might be too far away, try moving closer.")

# --- Insert a new paragraph right after this one (still before the trailing
#     bookmark), inheriting the same list/shading/spacing formatting, and
#     give it the new bullet text. ---
$para = $d.Paragraphs.Item($targetIndex)
$insertionPoint = $para.Range.End - 1
$ins = $d.Range($insertionPoint, $insertionPoint)
$ins.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newPara.Range.Text = "The Bluetooth isn’t paired. Press the button down on the Bluetooth to put into pairing mode. You should only need to do this once."

Write-Output ("Updated paragraph " + $targetIndex + ": " + $d.Paragraphs.Item($targetIndex).Range.Text)
Write-Output ("New paragraph " + ($targetIndex + 1) + ": " + $d.Paragraphs.Item($targetIndex + 1).Range.Text)
